$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 64.818184
$ws.Range("I42").Value = 36
$ws.Range("J42").Value = 81.28570999999999
$ws.Range("K42").Value = 108
$ws.Range("L42").Value = 243.85713
$ws.Range("M42").Value = 122
$ws.Range("N42").Value = -703.85713

# Row 106
$ws.Range("H106").Value = 3247.5
$ws.Range("I106").Value = 1995
$ws.Range("K106").Value = 1995
$ws.Range("M106").Value = -1364

# Row 116
$ws.Range("H116").Value = 2321.8572
$ws.Range("I116").Value = 2137.9333
$ws.Range("J116").Value = 2781.6667
$ws.Range("K116").Value = 2137.9333
$ws.Range("L116").Value = 2781.6667
$ws.Range("M116").Value = 1304.0667
$ws.Range("N116").Value = -9665.6667

# Row 129
$ws.Range("H129").Value = 1199.0483
$ws.Range("I129").Value = 355.69232
$ws.Range("J129").Value = 1422.7959
$ws.Range("K129").Value = 1067.07696
$ws.Range("L129").Value = 4268.3877
$ws.Range("M129").Value = 3932.92304
$ws.Range("N129").Value = -14268.3877

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 707.5
$ws.Range("I2").Value = 449.2
$ws.Range("K2").Value = 449.2
$ws.Range("M2").Value = -336.2

# Row 61
$ws.Range("H61").Value = 18521036
$ws.Range("I61").Value = 22729908
$ws.Range("K61").Value = 22729908
$ws.Range("M61").Value = -22729696

# Row 63
$ws.Range("H63").Value = 71429580
$ws.Range("I63").Value = 125001150
$ws.Range("J63").Value = 835.3333
$ws.Range("K63").Value = 125001150
$ws.Range("L63").Value = 835.3333
$ws.Range("M63").Value = -125000464
$ws.Range("N63").Value = -2207.3333

# Row 66
$ws.Range("H66").Value = 71429580
$ws.Range("I66").Value = 125001150
$ws.Range("J66").Value = 835.3333
$ws.Range("K66").Value = 625005750
$ws.Range("L66").Value = 4176.6665
$ws.Range("M66").Value = -625002318
$ws.Range("N66").Value = -11040.6665

# Row 76
$ws.Range("H76").Value = 20500
$ws.Range("J76").Value = 20500
$ws.Range("L76").Value = 20500
$ws.Range("N76").Value = -21176

# Row 79
$ws.Range("H79").Value = 20500
$ws.Range("J79").Value = 20500
$ws.Range("L79").Value = 20500
$ws.Range("N79").Value = -22840

# Row 88
$ws.Range("H88").Value = 9676.5
$ws.Range("I88").Value = 4853
$ws.Range("J88").Value = 14500
$ws.Range("K88").Value = 4853
$ws.Range("L88").Value = 14500
$ws.Range("M88").Value = -4447
$ws.Range("N88").Value = -15312

# Row 91
$ws.Range("H91").Value = 9676.5
$ws.Range("I91").Value = 4853
$ws.Range("J91").Value = 14500
$ws.Range("K91").Value = 4853
$ws.Range("L91").Value = 14500
$ws.Range("M91").Value = -3449
$ws.Range("N91").Value = -17308

# Row 116
$ws.Range("H116").Value = 707.5
$ws.Range("I116").Value = 449.2
$ws.Range("K116").Value = 449.2
$ws.Range("M116").Value = 1844.8

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 132
$ws.Range("H132").Value = 8622400
$ws.Range("I132").Value = 13159337
$ws.Range("J132").Value = 2218
$ws.Range("K132").Value = 39478011
$ws.Range("L132").Value = 6654
$ws.Range("M132").Value = -39475481
$ws.Range("N132").Value = -11714

# Row 136
$ws.Range("H136").Value = 18521036
$ws.Range("I136").Value = 22729908
$ws.Range("K136").Value = 68189724
$ws.Range("M136").Value = -68187174

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 707.5
$ws.Range("I3").Value = 449.2
$ws.Range("K3").Value = 449.2
$ws.Range("M3").Value = -335.2

# Row 86
$ws.Range("H86").Value = 27780488
$ws.Range("I86").Value = 2645.4546
$ws.Range("J86").Value = 71431384
$ws.Range("K86").Value = 2645.4546
$ws.Range("L86").Value = 71431384
$ws.Range("M86").Value = -1522.4546
$ws.Range("N86").Value = -71433630

# Row 89
$ws.Range("H89").Value = 27780488
$ws.Range("I89").Value = 2645.4546
$ws.Range("J89").Value = 71431384
$ws.Range("K89").Value = 13227.273
$ws.Range("L89").Value = 357156920
$ws.Range("M89").Value = -7611.273000000001
$ws.Range("N89").Value = -357168152

# Row 112
$ws.Range("H112").Value = 37549.285
$ws.Range("J112").Value = 37549.285
$ws.Range("L112").Value = 37549.285
$ws.Range("N112").Value = -40503.285

# Row 118
$ws.Range("H118").Value = 8166.6665
$ws.Range("J118").Value = 8166.6665
$ws.Range("L118").Value = 8166.6665
$ws.Range("N118").Value = -11480.6665

# Row 126
$ws.Range("H126").Value = 48000
$ws.Range("J126").Value = 48000
$ws.Range("L126").Value = 48000
$ws.Range("N126").Value = -57880

# Row 130
$ws.Range("H130").Value = 40780
$ws.Range("J130").Value = 40780
$ws.Range("L130").Value = 40780
$ws.Range("N130").Value = -50820

# Row 134
$ws.Range("H134").Value = 2905.1177
$ws.Range("I134").Value = 1792.6052
$ws.Range("J134").Value = 6157.077
$ws.Range("K134").Value = 5377.8156
$ws.Range("L134").Value = 18471.231
$ws.Range("M134").Value = -2842.8156
$ws.Range("N134").Value = -23541.231

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4978346.5
$ws.Range("I31").Value = 3541.52
$ws.Range("J31").Value = 19610126
$ws.Range("K31").Value = 3541.52
$ws.Range("L31").Value = 19610126
$ws.Range("M31").Value = -3246.52
$ws.Range("N31").Value = -19610716

# Row 34
$ws.Range("H34").Value = 4978346.5
$ws.Range("I34").Value = 3541.52
$ws.Range("J34").Value = 19610126
$ws.Range("K34").Value = 3541.52
$ws.Range("L34").Value = 19610126
$ws.Range("M34").Value = -3339.52
$ws.Range("N34").Value = -19610530

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 35.625
$ws.Range("I2").Value = 36.666668
$ws.Range("J2").Value = 32.5
$ws.Range("K2").Value = 36.666668
$ws.Range("L2").Value = 32.5
$ws.Range("M2").Value = 76.333332
$ws.Range("N2").Value = -258.5

# Row 3
$ws.Range("H3").Value = 455183.28
$ws.Range("I3").Value = 1250125
$ws.Range("J3").Value = 930.8570999999999
$ws.Range("K3").Value = 1250125
$ws.Range("L3").Value = 930.8570999999999
$ws.Range("M3").Value = -1250009
$ws.Range("N3").Value = -1162.8571

# Row 5
$ws.Range("H5").Value = 952
$ws.Range("I5").Value = 952
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 952
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -840
$ws.Range("N5").ClearContents()

# Row 7
$ws.Range("H7").Value = 2000000
$ws.Range("I7").Value = 2000000
$ws.Range("K7").Value = 2000000
$ws.Range("M7").Value = -1999888

# Row 8
$ws.Range("H8").Value = 2000000
$ws.Range("I8").Value = 2000000
$ws.Range("K8").Value = 2000000
$ws.Range("M8").Value = -1999861

# Row 9
$ws.Range("H9").Value = 150
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 150
$ws.Range("M9").Value = 20

# Row 10
$ws.Range("H10").Value = 5250640
$ws.Range("I10").Value = 5250640
$ws.Range("K10").Value = 5250640
$ws.Range("M10").Value = -5250471

# Row 80
$ws.Range("H80").Value = 14197547
$ws.Range("I80").Value = 33335774
$ws.Range("J80").Value = 2236154.2
$ws.Range("K80").Value = 33335774
$ws.Range("L80").Value = 2236154.2
$ws.Range("M80").Value = -33334776
$ws.Range("N80").Value = -2238150.2

# Row 83
$ws.Range("H83").Value = 14197547
$ws.Range("I83").Value = 33335774
$ws.Range("J83").Value = 2236154.2
$ws.Range("K83").Value = 166678870
$ws.Range("L83").Value = 11180771
$ws.Range("M83").Value = -166673878
$ws.Range("N83").Value = -11190755

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1929.5217
$ws.Range("I82").Value = 1482.3636
$ws.Range("J82").Value = 2339.4167
$ws.Range("K82").Value = 1482.3636
$ws.Range("L82").Value = 2339.4167
$ws.Range("M82").Value = -1121.3636
$ws.Range("N82").Value = -3061.4167

# Row 85
$ws.Range("H85").Value = 1929.5217
$ws.Range("I85").Value = 1482.3636
$ws.Range("J85").Value = 2339.4167
$ws.Range("K85").Value = 1482.3636
$ws.Range("L85").Value = 2339.4167
$ws.Range("M85").Value = -234.3635999999999
$ws.Range("N85").Value = -4835.4167

# Row 122
$ws.Range("H122").Value = 9525
$ws.Range("I122").Value = 15450
$ws.Range("J122").Value = 6562.5
$ws.Range("K122").Value = 46350
$ws.Range("L122").Value = 19687.5
$ws.Range("M122").Value = -43900
$ws.Range("N122").Value = -24587.5
